$wb = $excel.ActiveWorkbook

# Sheet 1: ALC
$ws = $wb.Worksheets.Item(1)
# Row 107
$ws.Range("H107").Value = 365.45456
$ws.Range("I107").Value = 392.23077
$ws.Range("J107").Value = 326.77777
$ws.Range("K107").Value = 392.23077
$ws.Range("L107").Value = 326.77777
$ws.Range("M107").Value = 1527.76923
$ws.Range("N107").Value = -4166.77777
# Row 112
$ws.Range("H112").Value = 1076.0454
$ws.Range("I112").Value = 666.9
$ws.Range("J112").Value = 1149.1072
$ws.Range("K112").Value = 2000.7
$ws.Range("L112").Value = 3447.3216
$ws.Range("M112").Value = -892.6999999999998
$ws.Range("N112").Value = -5663.321599999999
# Row 125
$ws.Range("H125").Value = 1181.7142
$ws.Range("I125").Value = 878
$ws.Range("J125").Value = 1586.6666
$ws.Range("K125").Value = 7902
$ws.Range("L125").Value = 14279.9994
$ws.Range("M125").Value = -5442
$ws.Range("N125").Value = -19199.9994

# Sheet 2: ARM
$ws = $wb.Worksheets.Item(2)
# Row 2
$ws.Range("H2").Value = 1102
$ws.Range("I2").Value = 1141.2222
$ws.Range("J2").Value = 866.6667
$ws.Range("K2").Value = 1141.2222
$ws.Range("L2").Value = 866.6667
$ws.Range("M2").Value = -1028.2222
$ws.Range("N2").Value = -1092.6667
# Row 61
$ws.Range("H61").Value = 1218.4736
$ws.Range("I61").Value = 1218.4736
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1218.4736
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -1006.4736
# Row 74
$ws.Range("H74").Value = 1277.1063
$ws.Range("I74").Value = 774.0476
$ws.Range("K74").Value = 774.0476
$ws.Range("M74").Value = 99.95240000000001
# Row 77
$ws.Range("H77").Value = 1277.1063
$ws.Range("I77").Value = 774.0476
$ws.Range("K77").Value = 3870.238
$ws.Range("M77").Value = 497.7620000000002
# Row 101
$ws.Range("H101").Value = 36726.25
$ws.Range("J101").Value = 36726.25
$ws.Range("L101").Value = 36726.25
$ws.Range("N101").Value = -43216.25
# Row 110
$ws.Range("I110").Value = 667.2727
$ws.Range("J110").Value = 1614.5
$ws.Range("K110").Value = 667.2727
$ws.Range("L110").Value = 1614.5
$ws.Range("M110").Value = 1377.7273
$ws.Range("N110").Value = -5704.5
# Row 111
$ws.Range("H111").Value = 30644
$ws.Range("J111").Value = 30644
$ws.Range("L111").Value = 30644
$ws.Range("N111").Value = -38824
# Row 116
$ws.Range("H116").Value = 1102
$ws.Range("I116").Value = 1141.2222
$ws.Range("J116").Value = 866.6667
$ws.Range("K116").Value = 1141.2222
$ws.Range("L116").Value = 866.6667
$ws.Range("M116").Value = 1152.7778
$ws.Range("N116").Value = -5454.6667
# Row 132
$ws.Range("H132").Value = 4207.381
$ws.Range("I132").Value = 4800.387
$ws.Range("K132").Value = 14401.161
$ws.Range("M132").Value = -11871.161
# Row 136
$ws.Range("H136").Value = 1218.4736
$ws.Range("I136").Value = 1218.4736
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 3655.4208
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -1105.4208

# Sheet 3: BSM
$ws = $wb.Worksheets.Item(3)
# Row 3
$ws.Range("H3").Value = 1102
$ws.Range("I3").Value = 1141.2222
$ws.Range("J3").Value = 866.6667
$ws.Range("K3").Value = 1141.2222
$ws.Range("L3").Value = 866.6667
$ws.Range("M3").Value = -1027.2222
$ws.Range("N3").Value = -1094.6667
# Row 100
$ws.Range("H100").Value = 25276.875
$ws.Range("J100").Value = 25276.875
$ws.Range("L100").Value = 25276.875
$ws.Range("N100").Value = -27440.875
# Row 107
$ws.Range("H107").Value = 820.1111
$ws.Range("I107").Value = 654
$ws.Range("J107").Value = 1027.75
$ws.Range("K107").Value = 654
$ws.Range("L107").Value = 1027.75
$ws.Range("M107").Value = 1266
$ws.Range("N107").Value = -4867.75
# Row 134
$ws.Range("H134").Value = 3968.658
$ws.Range("I134").Value = 4089.0881
$ws.Range("J134").Value = 2945
$ws.Range("K134").Value = 12267.2643
$ws.Range("L134").Value = 8835
$ws.Range("M134").Value = -9732.264299999999
$ws.Range("N134").Value = -13905

# Sheet 4: CRP
$ws = $wb.Worksheets.Item(4)
# Row 33
$ws.Range("H33").Value = 5176
$ws.Range("I33").Value = 265.5
$ws.Range("J33").Value = 14997
$ws.Range("K33").Value = 265.5
$ws.Range("L33").Value = 14997
$ws.Range("M33").Value = 113.5
$ws.Range("N33").Value = -15755
# Row 107
$ws.Range("H107").Value = 955.3043
$ws.Range("I107").Value = 1536.3334
$ws.Range("J107").Value = 321.45456
$ws.Range("K107").Value = 1536.3334
$ws.Range("L107").Value = 321.45456
$ws.Range("M107").Value = 383.6666
$ws.Range("N107").Value = -4161.45456
# Row 134
$ws.Range("H134").Value = 5030.5
$ws.Range("I134").Value = 1360.238
$ws.Range("J134").Value = 20445.6
$ws.Range("K134").Value = 4080.714
$ws.Range("L134").Value = 61336.8
$ws.Range("M134").Value = -1545.714
$ws.Range("N134").Value = -66406.79999999999

# Sheet 5: CUL
$ws = $wb.Worksheets.Item(5)
# Row 118
$ws.Range("H118").Value = 2426.2354
$ws.Range("I118").Value = 775
$ws.Range("K118").Value = 2325
$ws.Range("M118").Value = -1082

# Sheet 6: GSM
$ws = $wb.Worksheets.Item(6)
# Row 122
$ws.Range("H122").Value = 1728.4762
$ws.Range("I122").Value = 1450.625
$ws.Range("J122").Value = 2617.6
$ws.Range("K122").Value = 4351.875
$ws.Range("L122").Value = 7852.799999999999
$ws.Range("M122").Value = -1901.875
$ws.Range("N122").Value = -12752.8
# Row 132
$ws.Range("H132").Value = 4863.0586
$ws.Range("I132").Value = 5083.2856
$ws.Range("J132").Value = 3835.3333
$ws.Range("K132").Value = 15249.8568
$ws.Range("L132").Value = 11505.9999
$ws.Range("M132").Value = -12719.8568
$ws.Range("N132").Value = -16565.9999

# Sheet 7: LTW
$ws = $wb.Worksheets.Item(7)
# Row 40
$ws.Range("H40").Value = 2398
$ws.Range("I40").Value = 2197.5
$ws.Range("J40").Value = 3701.25
$ws.Range("K40").Value = 2197.5
$ws.Range("L40").Value = 3701.25
$ws.Range("M40").Value = -2061.5
$ws.Range("N40").Value = -3973.25
# Row 61
$ws.Range("H61").Value = 1535.7273
$ws.Range("I61").Value = 1434.2222
$ws.Range("J61").Value = 1992.5
$ws.Range("K61").Value = 1434.2222
$ws.Range("L61").Value = 1992.5
$ws.Range("M61").Value = -1232.2222
$ws.Range("N61").Value = -2396.5
# Row 68
$ws.Range("H68").Value = 2744.06
$ws.Range("I68").Value = 2137.5
$ws.Range("J68").Value = 2859.5952
$ws.Range("K68").Value = 2137.5
$ws.Range("L68").Value = 2859.5952
$ws.Range("M68").Value = -1388.5
$ws.Range("N68").Value = -4357.5952
# Row 71
$ws.Range("H71").Value = 2744.06
$ws.Range("I71").Value = 2137.5
$ws.Range("J71").Value = 2859.5952
$ws.Range("K71").Value = 10687.5
$ws.Range("L71").Value = 14297.976
$ws.Range("M71").Value = -6943.5
$ws.Range("N71").Value = -21785.976
# Row 113
$ws.Range("H113").Value = 1535.7273
$ws.Range("I113").Value = 1434.2222
$ws.Range("J113").Value = 1992.5
$ws.Range("K113").Value = 1434.2222
$ws.Range("L113").Value = 1992.5
$ws.Range("M113").Value = 735.7778000000001
$ws.Range("N113").Value = -6332.5
# Row 132
$ws.Range("H132").Value = 2789.6553
$ws.Range("I132").Value = 2187.5122
$ws.Range("J132").Value = 4241.8823
$ws.Range("K132").Value = 6562.5366
$ws.Range("L132").Value = 12725.6469
$ws.Range("M132").Value = -4032.5366
$ws.Range("N132").Value = -17785.6469
# Row 136
$ws.Range("H136").Value = 2187.389
$ws.Range("I136").Value = 1691.2084
$ws.Range("J136").Value = 3179.75
$ws.Range("K136").Value = 5073.6252
$ws.Range("L136").Value = 9539.25
$ws.Range("M136").Value = -2523.6252
$ws.Range("N136").Value = -14639.25

# Sheet 8: WVR
$ws = $wb.Worksheets.Item(8)
# Row 107
$ws.Range("H107").Value = 200
$ws.Range("I107").Value = 200
$ws.Range("K107").Value = 600
$ws.Range("M107").Value = 1320
# Row 122
$ws.Range("H122").Value = 2033729.1
$ws.Range("I122").Value = 1017003.8
$ws.Range("J122").Value = 5210995.5
$ws.Range("K122").Value = 3051011.4
$ws.Range("L122").Value = 15632986.5
$ws.Range("M122").Value = -3048561.4
$ws.Range("N122").Value = -15637886.5
# Row 132
$ws.Range("H132").Value = 2978.5095
$ws.Range("I132").Value = 3011.8333
$ws.Range("J132").Value = 2658.6
$ws.Range("K132").Value = 9035.499899999999
$ws.Range("L132").Value = 7975.799999999999
$ws.Range("M132").Value = -6505.499899999999
$ws.Range("N132").Value = -13035.8
# Row 136
$ws.Range("H136").Value = 847.9655
$ws.Range("I136").Value = 803.37036
$ws.Range("J136").Value = 1450
$ws.Range("K136").Value = 2410.11108
$ws.Range("L136").Value = 4350
$ws.Range("M136").Value = 139.8889199999999
$ws.Range("N136").Value = -9450
